$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (COM ColumnWidth=30.17 round-trips to a stored char width of 31)
$ws.Columns.Item(1).ColumnWidth = 30.17

# New data rows for A2:D17. Row 2 previously held PAMPERS-PANTS-PC-MES-Gx64,
# which now moves down to row 9 among many newly inserted rows.
$data = @(
    @("PAMPERS-PREMIUM-MES-HP-XGx58", 4895, 5482, 10),
    @("BABYSEC-ULTRA-REGULAR-Gx8", 451, 505, 5),
    @("BABYSEC-ULTRA-REGULAR-Mx8", 451, 505, 5),
    @("BABYSEC-ULTRA-REGULAR-XXGx8", 451, 505, 5),
    @("BABYSEC-ULTRA-REGULAR-Px12", 451, 505, 5),
    @("BABYSEC-ULTRA-REGULAR-XGx8", 451, 505, 5),
    @("PAMPERS-PREMIUM-XTR-RN-x36", 1950, 2184, 5),
    @("PAMPERS-PANTS-PC-MES-Gx64", 5201, 5825, 5),
    @("PAMPERS-CONFORT-XTR-MES-XXGx54", 4405, 4934, 5),
    @("PAMPERS-PREMIUM-XTR-MES-XXGx54", 4895, 5482, 5),
    @("PAMPERS-CONFORT-XTR-HIPER-Px56", 2501, 2801, 5),
    @("PAMPERS-SUPER-XTR-REG-Gx9", 526, 589, 5),
    @("BABYSEC-PREMIUM-JUMBO-XGx48", 3357, 3760, 4),
    @("BABYSEC-PREMIUM-JUMBO-XXGx44", 3357, 3760, 4),
    @("BABYSEC-PREMIUM-JUMBO-Gx60", 3357, 3760, 4),
    @("PAMPERS-TOTAL-PROTECT-XXGx34", 2110, 2363, 2)
)

# Pre-format the new rows (3:17) so every new cell carries the same centered
# style (s="2") that row 2 already uses, before the values are written in.
$ws.Range("A3:D17").HorizontalAlignment = -4108  # xlCenter

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}
